$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 8179.385
$ws.Range("I9").Value = 347.57144
$ws.Range("K9").Value = 347.57144
$ws.Range("M9").Value = -178.57144
$ws.Range("H17").Value = 898.64703
$ws.Range("J17").Value = 917.3125
$ws.Range("L17").Value = 2751.9375
$ws.Range("N17").Value = -3087.9375
$ws.Range("H32").Value = 4124.25
$ws.Range("I32").Value = 3299.25
$ws.Range("J32").Value = 4536.75
$ws.Range("K32").Value = 3299.25
$ws.Range("L32").Value = 4536.75
$ws.Range("M32").Value = -2973.25
$ws.Range("N32").Value = -5188.75
$ws.Range("H74").Value = 18363.334
$ws.Range("I74").Value = 5090
$ws.Range("K74").Value = 5090
$ws.Range("M74").Value = -4154
$ws.Range("H76").Value = 11263.917
$ws.Range("I76").Value = 11067.5
$ws.Range("K76").Value = 11067.5
$ws.Range("M76").Value = -10752.5
$ws.Range("H77").Value = 18363.334
$ws.Range("I77").Value = 5090
$ws.Range("K77").Value = 25450
$ws.Range("M77").Value = -20770
$ws.Range("H79").Value = 11263.917
$ws.Range("I79").Value = 11067.5
$ws.Range("K79").Value = 11067.5
$ws.Range("M79").Value = -9975.5
$ws.Range("H92").Value = 4957.25
$ws.Range("I92").Value = 4967
$ws.Range("J92").Value = 4947.5
$ws.Range("K92").Value = 4967
$ws.Range("L92").Value = 4947.5
$ws.Range("M92").Value = -3719
$ws.Range("N92").Value = -7443.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2705.818
$ws.Range("J2").Value = 3003.3333
$ws.Range("L2").Value = 3003.3333
$ws.Range("N2").Value = -3229.3333
$ws.Range("H32").Value = 6273.119
$ws.Range("I32").Value = 5173.2256
$ws.Range("K32").Value = 5173.2256
$ws.Range("M32").Value = -4886.2256
$ws.Range("H61").Value = 3155.6057
$ws.Range("I61").Value = 2337.0667
$ws.Range("K61").Value = 2337.0667
$ws.Range("M61").Value = -2125.0667
$ws.Range("H74").Value = 37229.66
$ws.Range("I74").Value = 22621.523
$ws.Range("K74").Value = 22621.523
$ws.Range("M74").Value = -21747.523
$ws.Range("H77").Value = 37229.66
$ws.Range("I77").Value = 22621.523
$ws.Range("K77").Value = 113107.615
$ws.Range("M77").Value = -108739.615
$ws.Range("H116").Value = 2705.818
$ws.Range("J116").Value = 3003.3333
$ws.Range("L116").Value = 3003.3333
$ws.Range("N116").Value = -7591.3333
$ws.Range("H132").Value = 1839.5454
$ws.Range("I132").Value = 1886.0312
$ws.Range("J132").Value = 1715.5834
$ws.Range("K132").Value = 5658.0936
$ws.Range("L132").Value = 5146.7502
$ws.Range("M132").Value = -3128.0936
$ws.Range("N132").Value = -10206.7502
$ws.Range("H136").Value = 3155.6057
$ws.Range("I136").Value = 2337.0667
$ws.Range("K136").Value = 7011.2001
$ws.Range("M136").Value = -4461.2001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2705.818
$ws.Range("J3").Value = 3003.3333
$ws.Range("L3").Value = 3003.3333
$ws.Range("N3").Value = -3231.3333
$ws.Range("H86").Value = 62503040
$ws.Range("J86").Value = 3300
$ws.Range("L86").Value = 3300
$ws.Range("N86").Value = -5546
$ws.Range("H89").Value = 62503040
$ws.Range("J89").Value = 3300
$ws.Range("L89").Value = 16500
$ws.Range("N89").Value = -27732
$ws.Range("H105").Value = 125002776
$ws.Range("I105").Value = 125002776
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 125002776
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -125001029
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 18310.125
$ws.Range("I107").Value = 18310.125
$ws.Range("K107").Value = 18310.125
$ws.Range("M107").Value = -16390.125
$ws.Range("H134").Value = 1869.2041
$ws.Range("I134").Value = 1575.4147
$ws.Range("K134").Value = 4726.2441
$ws.Range("M134").Value = -2191.2441

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 11999
$ws.Range("J14").Value = 19999
$ws.Range("L14").Value = 19999
$ws.Range("N14").Value = -20339
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("H25").Value = 4054
$ws.Range("I25").Value = 1099.3334
$ws.Range("K25").Value = 1099.3334
$ws.Range("M25").Value = -925.3334
$ws.Range("H99").Value = 2756.568
$ws.Range("J99").Value = 2771.5715
$ws.Range("L99").Value = 2771.5715
$ws.Range("N99").Value = -5767.5715
$ws.Range("H126").Value = 2756.568
$ws.Range("J126").Value = 2771.5715
$ws.Range("L126").Value = 8314.7145
$ws.Range("N126").Value = -13254.7145
$ws.Range("H132").Value = 2898.7585
$ws.Range("I132").Value = 2372.7407
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 7118.222099999999
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -4588.222099999999
$ws.Range("N132").Value = -35060
$ws.Range("H134").Value = 2706.689
$ws.Range("I134").Value = 2734.1135
$ws.Range("K134").Value = 8202.3405
$ws.Range("M134").Value = -5667.3405

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 11500
$ws.Range("J64").Value = 12250
$ws.Range("L64").Value = 36750
$ws.Range("N64").Value = -37290
$ws.Range("H67").Value = 11500
$ws.Range("J67").Value = 12250
$ws.Range("L67").Value = 36750
$ws.Range("N67").Value = -38622
$ws.Range("H81").Value = 6897.5
$ws.Range("J81").Value = 8900
$ws.Range("L81").Value = 26700
$ws.Range("N81").Value = -28946
$ws.Range("H84").Value = 6897.5
$ws.Range("J84").Value = 8900
$ws.Range("L84").Value = 80100
$ws.Range("N84").Value = -91332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1711.826
$ws.Range("I80").Value = 1724.5555
$ws.Range("J80").Value = 1666
$ws.Range("K80").Value = 1724.5555
$ws.Range("L80").Value = 1666
$ws.Range("M80").Value = -726.5554999999999
$ws.Range("N80").Value = -3662
$ws.Range("H83").Value = 1711.826
$ws.Range("I83").Value = 1724.5555
$ws.Range("J83").Value = 1666
$ws.Range("K83").Value = 8622.7775
$ws.Range("L83").Value = 8330
$ws.Range("M83").Value = -3630.7775
$ws.Range("N83").Value = -18314

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 30198.2
$ws.Range("I43").Value = 28997.75
$ws.Range("K43").Value = 28997.75
$ws.Range("M43").Value = -28804.75
$ws.Range("H82").Value = 751.37933
$ws.Range("J82").Value = 963.3
$ws.Range("L82").Value = 963.3
$ws.Range("N82").Value = -1685.3
$ws.Range("H85").Value = 751.37933
$ws.Range("J85").Value = 963.3
$ws.Range("L85").Value = 963.3
$ws.Range("N85").Value = -3459.3
$ws.Range("H136").Value = 21933.244
$ws.Range("I136").Value = 1687.7693
$ws.Range("K136").Value = 5063.3079
$ws.Range("M136").Value = -2513.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16199.2
$ws.Range("I62").Value = 18623.25
$ws.Range("J62").Value = 14583.167
$ws.Range("K62").Value = 18623.25
$ws.Range("L62").Value = 14583.167
$ws.Range("M62").Value = -17999.25
$ws.Range("N62").Value = -15831.167
$ws.Range("H65").Value = 16199.2
$ws.Range("I65").Value = 18623.25
$ws.Range("J65").Value = 14583.167
$ws.Range("K65").Value = 93116.25
$ws.Range("L65").Value = 72915.83499999999
$ws.Range("M65").Value = -89996.25
$ws.Range("N65").Value = -79155.83499999999
$ws.Range("H132").Value = 3225.7083
$ws.Range("I132").Value = 3020.9
$ws.Range("J132").Value = 4249.75
$ws.Range("K132").Value = 9062.700000000001
$ws.Range("L132").Value = 12749.25
$ws.Range("M132").Value = -6532.700000000001
$ws.Range("N132").Value = -17809.25
